# Script auto-generated to apply the "Atualizado por script em 20-11-2023 08:45" edit.
# 1) Re-orders the F:V (match) data block among a set of existing rows that share
#    the same kickoff date (columns A-E stay put; only F..V move).
# 2) Appends 6 brand-new match rows (131-136) at the end and grows the dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Part 1: reshuffle F:V between rows (see $mapping: targetRow = sourceRow)
# ---------------------------------------------------------------------------

$mapping = [ordered]@{
    15 = 18
    17 = 15
    18 = 17
    24 = 25
    25 = 24
    31 = 33
    32 = 31
    33 = 32
    35 = 38
    36 = 37
    37 = 36
    38 = 35
    43 = 45
    45 = 43
    64 = 66
    65 = 64
    66 = 65
    67 = 68
    68 = 69
    69 = 67
    76 = 77
    77 = 76
    78 = 79
    79 = 78
    80 = 81
    81 = 82
    82 = 80
    84 = 87
    85 = 86
    86 = 85
    87 = 84
    93 = 95
    94 = 97
    95 = 96
    96 = 94
    97 = 93
    99 = 100
    100 = 99
    102 = 103
    103 = 102
    104 = 107
    105 = 106
    106 = 105
    107 = 104
    113 = 114
    114 = 113
    119 = 120
    120 = 119
}

# columns F(6) .. V(22); numeric columns vs text columns
$numericCols = @(7, 9, 10, 12, 14, 16, 18, 20)   # G, I, J, L, N, P, R, T
$firstCol = 6   # F
$lastCol  = 22  # V

# Snapshot every row that is used as a *source* before any writes happen,
# so a row that is both a source and a target elsewhere doesn't get
# clobbered before it's been read.
$snapshot = @{}
foreach ($srcRow in ($mapping.Values | Sort-Object -Unique)) {
    $rowVals = @{}
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $rowVals[$col] = $ws.Cells.Item($srcRow, $col).Value()
    }
    $snapshot[$srcRow] = $rowVals
}

foreach ($targetRow in $mapping.Keys) {
    $srcRow = $mapping[$targetRow]
    $rowVals = $snapshot[$srcRow]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $ws.Cells.Item($targetRow, $col).Value = $rowVals[$col]
    }
}

# ---------------------------------------------------------------------------
# Part 2: append the 6 new rows (131-136) at the bottom of the sheet
# ---------------------------------------------------------------------------

$newRows = @(
    ,@(131, 130, "italy", "serie-c-group-b", "2023-2024", 45248.67708333334, "Pineto", "1", "Perugia", "1", "3.84", "16/11/2023 09:12", "4.81", "18/11/2023 15:53", "3.12", "16/11/2023 09:12", "3.31", "18/11/2023 16:05", "1.92", "16/11/2023 09:12", "1.83", "18/11/2023 16:05", "https://www.betexplorer.com/football/italy/serie-c-group-b/pineto-perugia/hEFPJP8E/")
    ,@(132, 131, "italy", "serie-c-group-b", "2023-2024", 45248.77083333334, "Ancona", "2", "Recanatese", "0", "1.99", "16/11/2023 09:12", "2.06", "18/11/2023 18:28", "3.11", "16/11/2023 09:12", "3.26", "18/11/2023 18:28", "3.6", "16/11/2023 09:12", "3.76", "18/11/2023 18:28", "https://www.betexplorer.com/football/italy/serie-c-group-b/ancona-recanatese/pjJTa4Hr/")
    ,@(133, 132, "italy", "serie-c-group-b", "2023-2024", 45249.58333333334, "Fermana", "2", "Arezzo", "3", "3.01", "16/11/2023 09:12", "3.39", "19/11/2023 13:39", "2.85", "16/11/2023 09:12", "2.85", "19/11/2023 13:39", "2.39", "16/11/2023 09:12", "2.43", "19/11/2023 13:39", "https://www.betexplorer.com/football/italy/serie-c-group-b/fermana-arezzo/zuCtc211/")
    ,@(134, 133, "italy", "serie-c-group-b", "2023-2024", 45249.67708333334, "Cesena", "3", "Lucchese", "0", "1.48", "16/11/2023 21:12", "1.55", "19/11/2023 16:10", "3.76", "16/11/2023 21:12", "3.63", "19/11/2023 16:10", "6.19", "16/11/2023 21:12", "7.29", "19/11/2023 16:10", "https://www.betexplorer.com/football/italy/serie-c-group-b/cesena-lucchese/tICxbrne/")
    ,@(135, 134, "italy", "serie-c-group-b", "2023-2024", 45249.67708333334, "Carrarese", "1", "Spal", "0", "1.89", "16/11/2023 09:12", "1.66", "19/11/2023 16:12", "3.15", "16/11/2023 09:12", "3.43", "19/11/2023 16:12", "3.89", "16/11/2023 09:12", "6.12", "19/11/2023 16:12", "https://www.betexplorer.com/football/italy/serie-c-group-b/carrarese-spal/SE8YbOWl/")
    ,@(136, 135, "italy", "serie-c-group-b", "2023-2024", 45249.67708333334, "Gubbio", "5", "Sestri Levante", "2", "1.53", "16/11/2023 18:12", "1.6", "19/11/2023 16:11", "3.75", "16/11/2023 18:12", "3.61", "19/11/2023 16:11", "5.85", "16/11/2023 18:12", "6.56", "19/11/2023 16:11", "https://www.betexplorer.com/football/italy/serie-c-group-b/gubbio-sestri-levante/EVApdMG7/")
)

$rowNumericCols = @(2, 7, 9, 10, 12, 14, 16, 18, 20)  # A(Indice), G, I, J, L, N, P, R, T -- col E handled separately (always numeric)

foreach ($rowData in $newRows) {
    $rowNum = $rowData[0]
    for ($i = 1; $i -lt $rowData.Count; $i++) {
        $col = $i + 1   # data index 1 -> column B(2) ... but A is data index1? adjust below
    }

    # Explicit column writes (clearer & safer than the generic loop above)
    $ws.Cells.Item($rowNum, 1).Value  = [double]$rowData[1]     # A Indice
    $ws.Cells.Item($rowNum, 2).Value  = $rowData[2]             # B pais
    $ws.Cells.Item($rowNum, 3).Value  = $rowData[3]             # C torneio
    $ws.Cells.Item($rowNum, 4).Value  = $rowData[4]             # D temporada
    $ws.Cells.Item($rowNum, 5).Value  = [double]$rowData[5]     # E data_partida (serial)
    $ws.Cells.Item($rowNum, 6).Value  = $rowData[6]             # F home
    $ws.Cells.Item($rowNum, 7).Value  = [double]$rowData[7]     # G home_ft_gols
    $ws.Cells.Item($rowNum, 8).Value  = $rowData[8]             # H away
    $ws.Cells.Item($rowNum, 9).Value  = [double]$rowData[9]     # I away_ft_gols
    $ws.Cells.Item($rowNum, 10).Value = [double]$rowData[10]    # J home_opening_odds
    $ws.Cells.Item($rowNum, 11).Value = $rowData[11]            # K home_opening_data_hora
    $ws.Cells.Item($rowNum, 12).Value = [double]$rowData[12]    # L home_closing_odds
    $ws.Cells.Item($rowNum, 13).Value = $rowData[13]            # M home_closing_data_hora
    $ws.Cells.Item($rowNum, 14).Value = [double]$rowData[14]    # N draw_opening_odds
    $ws.Cells.Item($rowNum, 15).Value = $rowData[15]            # O draw_opening_data_hora
    $ws.Cells.Item($rowNum, 16).Value = [double]$rowData[16]    # P draw_closing_odds
    $ws.Cells.Item($rowNum, 17).Value = $rowData[17]            # Q draw_closing_data_hora
    $ws.Cells.Item($rowNum, 18).Value = [double]$rowData[18]    # R away_opening_odds
    $ws.Cells.Item($rowNum, 19).Value = $rowData[19]            # S away_opening_data_hora
    $ws.Cells.Item($rowNum, 20).Value = [double]$rowData[20]    # T away_closing_odds
    $ws.Cells.Item($rowNum, 21).Value = $rowData[21]            # U away_closing_data_hora
    $ws.Cells.Item($rowNum, 22).Value = $rowData[22]            # V url_partida
}

Write-Output "Edit complete."
